# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values. Recalculated K values replace the old Strike# values.
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 0
    10 = 2
    11 = 3
    12 = 0
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
